$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 292
$ws.Range("I5").Value = 240.25
$ws.Range("K5").Value = 240.25
$ws.Range("M5").Value = -125.25

$ws.Range("H11").Value = 122.94444
$ws.Range("I11").Value = 122.94444
$ws.Range("K11").Value = 122.94444
$ws.Range("M11").Value = 17.05556

$ws.Range("H40").Value = 1326.836
$ws.Range("I40").Value = 1133.8636
$ws.Range("J40").Value = 1435.6923
$ws.Range("K40").Value = 1133.8636
$ws.Range("L40").Value = 1435.6923
$ws.Range("M40").Value = -958.8635999999999
$ws.Range("N40").Value = -1785.6923

$ws.Range("H41").Value = 1588.1818
$ws.Range("I41").Value = 2235.8572
$ws.Range("J41").Value = 454.75
$ws.Range("K41").Value = 2235.8572
$ws.Range("L41").Value = 454.75
$ws.Range("M41").Value = -1795.8572
$ws.Range("N41").Value = -1334.75

$ws.Range("H53").Value = 194.46153
$ws.Range("I53").Value = 187
$ws.Range("J53").Value = 206.4
$ws.Range("K53").Value = 187
$ws.Range("L53").Value = 206.4
$ws.Range("M53").Value = 450
$ws.Range("N53").Value = -1480.4

$ws.Range("H70").Value = 4566.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4566.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 13699.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -14239.5

$ws.Range("H73").Value = 4566.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4566.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 13699.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -15571.5

$ws.Range("H98").Value = 864.3158
$ws.Range("I98").Value = 852.0714
$ws.Range("J98").Value = 898.6
$ws.Range("K98").Value = 852.0714
$ws.Range("L98").Value = 898.6
$ws.Range("M98").Value = 645.9286
$ws.Range("N98").Value = -3894.6

$ws.Range("H112").Value = 92595.27
$ws.Range("J112").Value = 101705.7
$ws.Range("L112").Value = 305117.1
$ws.Range("N112").Value = -307333.1

$ws.Range("H122").Value = 864.3158
$ws.Range("I122").Value = 852.0714
$ws.Range("J122").Value = 898.6
$ws.Range("K122").Value = 2556.2142
$ws.Range("L122").Value = 2695.8
$ws.Range("M122").Value = -106.2142000000003
$ws.Range("N122").Value = -7595.8

$ws.Range("H132").Value = 62201.516
$ws.Range("I132").Value = 40025
$ws.Range("K132").Value = 120075
$ws.Range("M132").Value = -117545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 3535.2
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H43").Value = 44751
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H45").Value = 8896.5
$ws.Range("I45").Value = 10418
$ws.Range("K45").Value = 10418
$ws.Range("M45").Value = -10041

$ws.Range("H61").Value = 4995.1
$ws.Range("I61").Value = 4539.8
$ws.Range("K61").Value = 4539.8
$ws.Range("M61").Value = -4327.8

$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496

$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716

$ws.Range("H74").Value = 8939.637000000001
$ws.Range("I74").Value = 1190.1666
$ws.Range("J74").Value = 18239
$ws.Range("K74").Value = 1190.1666
$ws.Range("L74").Value = 18239
$ws.Range("M74").Value = -316.1666
$ws.Range("N74").Value = -19987

$ws.Range("H77").Value = 8939.637000000001
$ws.Range("I77").Value = 1190.1666
$ws.Range("J77").Value = 18239
$ws.Range("K77").Value = 5950.833000000001
$ws.Range("L77").Value = 91195
$ws.Range("M77").Value = -1582.833000000001
$ws.Range("N77").Value = -99931

$ws.Range("H99").Value = 3535.2
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H122").Value = 2382.6428
$ws.Range("I122").Value = 2027.6666
$ws.Range("J122").Value = 3021.6
$ws.Range("K122").Value = 6082.9998
$ws.Range("L122").Value = 9064.799999999999
$ws.Range("M122").Value = -3632.9998
$ws.Range("N122").Value = -13964.8

$ws.Range("H132").Value = 2690.96
$ws.Range("I132").Value = 2212.85
$ws.Range("K132").Value = 6638.549999999999
$ws.Range("M132").Value = -4108.549999999999

$ws.Range("H135").Value = 90149.836
$ws.Range("J135").Value = 90149.836
$ws.Range("L135").Value = 90149.836
$ws.Range("N135").Value = -100289.836

$ws.Range("H136").Value = 4995.1
$ws.Range("I136").Value = 4539.8
$ws.Range("K136").Value = 13619.4
$ws.Range("M136").Value = -11069.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1847.0555
$ws.Range("J94").Value = 1276.3334
$ws.Range("L94").Value = 1276.3334
$ws.Range("N94").Value = -2178.3334

$ws.Range("H134").Value = 2303.111
$ws.Range("I134").Value = 2026.5161
$ws.Range("J134").Value = 4018
$ws.Range("K134").Value = 6079.5483
$ws.Range("L134").Value = 12054
$ws.Range("M134").Value = -3544.5483
$ws.Range("N134").Value = -17124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2815.9
$ws.Range("I3").Value = 2303.25
$ws.Range("J3").Value = 4866.5
$ws.Range("K3").Value = 2303.25
$ws.Range("L3").Value = 4866.5
$ws.Range("M3").Value = -2190.25
$ws.Range("N3").Value = -5092.5

$ws.Range("H26").Value = 9019
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H69").Value = 39833.168
$ws.Range("I69").Value = 39833.168
$ws.Range("K69").Value = 39833.168
$ws.Range("M69").Value = -39084.168

$ws.Range("H72").Value = 39833.168
$ws.Range("I72").Value = 39833.168
$ws.Range("K72").Value = 119499.504
$ws.Range("M72").Value = -115755.504

$ws.Range("H107").Value = 497.875
$ws.Range("I107").Value = 497.875
$ws.Range("K107").Value = 497.875
$ws.Range("M107").Value = 1422.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 225
$ws.Range("I25").Value = 225
$ws.Range("K25").Value = 675
$ws.Range("M25").Value = -506

$ws.Range("H30").Value = 225
$ws.Range("I30").Value = 225
$ws.Range("K30").Value = 675
$ws.Range("M30").Value = -573

$ws.Range("H40").Value = 321.5
$ws.Range("I40").Value = 209
$ws.Range("J40").Value = 434
$ws.Range("K40").Value = 836
$ws.Range("L40").Value = 1736
$ws.Range("M40").Value = -767
$ws.Range("N40").Value = -1874

$ws.Range("H132").Value = 1024.6666
$ws.Range("I132").Value = 1027.091
$ws.Range("J132").Value = 998
$ws.Range("K132").Value = 9243.819
$ws.Range("L132").Value = 8982
$ws.Range("M132").Value = -6713.819
$ws.Range("N132").Value = -14042

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 192464.17
$ws.Range("J42").Value = 188696.25
$ws.Range("L42").Value = 188696.25
$ws.Range("N42").Value = -189666.25

$ws.Range("H70").Value = 9952.056
$ws.Range("I70").Value = 8473
$ws.Range("J70").Value = 11135.3
$ws.Range("K70").Value = 8473
$ws.Range("L70").Value = 11135.3
$ws.Range("M70").Value = -8203
$ws.Range("N70").Value = -11675.3

$ws.Range("H73").Value = 9952.056
$ws.Range("I73").Value = 8473
$ws.Range("J73").Value = 11135.3
$ws.Range("K73").Value = 8473
$ws.Range("L73").Value = 11135.3
$ws.Range("M73").Value = -7537
$ws.Range("N73").Value = -13007.3

$ws.Range("H102").Value = 4890
$ws.Range("I102").Value = 4237.5
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 4237.5
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -2615.5
$ws.Range("N102").Value = -10744

$ws.Range("H115").Value = 192464.17
$ws.Range("J115").Value = 188696.25
$ws.Range("L115").Value = 188696.25
$ws.Range("N115").Value = -191046.25

$ws.Range("H126").Value = 43250.625
$ws.Range("I126").Value = 48715
$ws.Range("K126").Value = 146145
$ws.Range("M126").Value = -143675

$ws.Range("H132").Value = 373533.03
$ws.Range("I132").Value = 479328.2
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 1437984.6
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -1435454.6
$ws.Range("N132").Value = -14810

$ws.Range("H135").Value = 55713.332
$ws.Range("J135").Value = 55713.332
$ws.Range("L135").Value = 55713.332
$ws.Range("N135").Value = -65853.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 310.75
$ws.Range("I55").Value = 320.63635
$ws.Range("J55").Value = 289
$ws.Range("K55").Value = 320.63635
$ws.Range("L55").Value = 289
$ws.Range("M55").Value = -147.63635
$ws.Range("N55").Value = -635

$ws.Range("H108").Value = 39975
$ws.Range("J108").Value = 39975
$ws.Range("L108").Value = 39975
$ws.Range("N108").Value = -47655
